# Usman - New users added to cobaltUsers.xls
# Adds 33 new user rows (81-113) to the "Users" sheet: KHPaddUser1-6,
# SearchWhatsMarketUser1-8, SearchKnowHowUser1-8, AskUser1-6, AssetPageUser1-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Row data: Name goes in column A, Email (hyperlinked via mailto:) in column G.
# Columns B/E/F repeat the same constant values used by every row in this block.
$rows = @(
    @{Row=81; Name='KHPaddUser1'; Email='KHPaddUser1@mailinator.com'},
    @{Row=82; Name='KHPaddUser2'; Email='KHPaddUser2@mailinator.com'},
    @{Row=83; Name='KHPaddUser3'; Email='KHPaddUser3@mailinator.com'},
    @{Row=84; Name='KHPaddUser4'; Email='KHPaddUser4@mailinator.com'},
    @{Row=85; Name='KHPaddUser5'; Email='KHPaddUser5@mailinator.com'},
    @{Row=86; Name='KHPaddUser6'; Email='KHPaddUser6@mailinator.com'},
    @{Row=87; Name='SearchWhatsMarketUser1'; Email='SearchWhatsMarketUser1@mailinator.com '},
    @{Row=88; Name='SearchWhatsMarketUser2'; Email='SearchWhatsMarketUser2@mailinator.com '},
    @{Row=89; Name='SearchWhatsMarketUser3'; Email='SearchWhatsMarketUser3@mailinator.com '},
    @{Row=90; Name='SearchWhatsMarketUser4'; Email='SearchWhatsMarketUser4@mailinator.com '},
    @{Row=91; Name='SearchWhatsMarketUser5'; Email='SearchWhatsMarketUser5@mailinator.com '},
    @{Row=92; Name='SearchWhatsMarketUser6'; Email='SearchWhatsMarketUser6@mailinator.com '},
    @{Row=93; Name='SearchWhatsMarketUser7'; Email='SearchWhatsMarketUser7@mailinator.com '},
    @{Row=94; Name='SearchWhatsMarketUser8'; Email='SearchWhatsMarketUser8@mailinator.com '},
    @{Row=95; Name='SearchKnowHowUser1'; Email='SearchKnowHowUser1@mailinator.com '},
    @{Row=96; Name='SearchKnowHowUser2'; Email='SearchKnowHowUser2@mailinator.com '},
    @{Row=97; Name='SearchKnowHowUser3'; Email='SearchKnowHowUser3@mailinator.com '},
    @{Row=98; Name='SearchKnowHowUser4'; Email='SearchKnowHowUser4@mailinator.com '},
    @{Row=99; Name='SearchKnowHowUser5'; Email='SearchKnowHowUser5@mailinator.com '},
    @{Row=100; Name='SearchKnowHowUser6'; Email='SearchKnowHowUser6@mailinator.com '},
    @{Row=101; Name='SearchKnowHowUser7'; Email='SearchKnowHowUser7@mailinator.com '},
    @{Row=102; Name='SearchKnowHowUser8'; Email='SearchKnowHowUser8@mailinator.com '},
    @{Row=103; Name='AskUser1'; Email='AskUser1@mailinator.com '},
    @{Row=104; Name='AskUser2'; Email='AskUser2@mailinator.com '},
    @{Row=105; Name='AskUser3'; Email='AskUser3@mailinator.com '},
    @{Row=106; Name='AskUser4'; Email='AskUser4@mailinator.com '},
    @{Row=107; Name='AskUser5'; Email='AskUser5@mailinator.com '},
    @{Row=108; Name='AskUser6'; Email='AskUser6@mailinator.com '},
    @{Row=109; Name='AssetPageUser1'; Email='AssetPageUser1@mailinator.com '},
    @{Row=110; Name='AssetPageUser2'; Email='AssetPageUser2@mailinator.com '},
    @{Row=111; Name='AssetPageUser3'; Email='AssetPageUser3@mailinator.com '},
    @{Row=112; Name='AssetPageUser4'; Email='AssetPageUser4@mailinator.com '},
    @{Row=113; Name='AssetPageUser5'; Email='AssetPageUser5@mailinator.com '}
)
$byRow = @{}
foreach ($r in $rows) { $byRow[$r.Row] = $r }

# Rows that (per the source workbook) never got their mailto hyperlink wired
# up -- the email text is present but is a plain string, not a hyperlink.
$noHyperlink = @(81, 112)

# Hyperlinks are wired up in this exact (non-sequential) order in the source
# file -- row 113's link was added before row 111's, which is also the one
# row whose link kept a stale display caption from copy/paste (it shows
# AssetPageUser1's address even though the cell/link point at
# AssetPageUser3).
$hyperlinkOrder = @(82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,113,111)
$staleDisplay = @{111 = 'AssetPageUser1@mailinator.com '}

# First lay down all the plain cell values/text in row order.
foreach ($r in $rows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.Name
    $ws.Range("B$row").Value = "Password1"
    $ws.Range("E$row").Value = "THIS IS IN USE 24/7 - DO NOT USE!"
    $ws.Range("F$row").Value = "N"
    $ws.Range("G$row").Value = $r.Email

    # Re-apply the E:G number format/style used throughout this block (copied
    # from the row directly above the new block, row 80).
    $ws.Range("E80:G80").Copy() | Out-Null
    $ws.Range("E$row`:G$row").PasteSpecial(-4122) | Out-Null

    # Column A in this new block uses a dedicated Arial 10 font (distinct
    # from the default font used by the rest of the sheet).
    $ws.Range("A$row").Font.Name = "Arial"
    $ws.Range("A$row").Font.Size = 10
}

# Now wire up the mailto hyperlinks in their recorded order.
foreach ($row in $hyperlinkOrder) {
    $r = $byRow[$row]
    $target = "mailto:" + $r.Email.Trim()

    if ($staleDisplay.ContainsKey($row)) {
        $ws.Hyperlinks.Add($ws.Range("G$row"), $target, "", "", $staleDisplay[$row]) | Out-Null
        # Hyperlinks.Add drives the cell text from the display caption above;
        # force the cell back to this row's own email so only the link
        # metadata keeps the stale caption.
        $ws.Range("G$row").Value = $r.Email
    } else {
        $ws.Hyperlinks.Add($ws.Range("G$row"), $target) | Out-Null
    }

    # Hyperlinks.Add re-applies its own "Hyperlink" cell style; restore the
    # E:G format used by the rest of this block.
    $ws.Range("E80:G80").Copy() | Out-Null
    $ws.Range("E$row`:G$row").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false

# Match the saved selection/active-cell state recorded in the edited workbook.
$ws.Range("A81:G113").Select() | Out-Null
